$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F12").Value = -5
$ws.Range("F16").Value = -4
$ws.Range("F19").Value = -2
$ws.Range("F22").Value = 3
$ws.Range("F25").Value = 5
$ws.Range("F26").Value = 2
$ws.Range("F29").Value = 0
$ws.Range("F32").Value = -3
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = -2
$ws.Range("F36").Value = 2
$ws.Range("F38").Value = 2
$ws.Range("F43").Value = 4
$ws.Range("F47").Value = -2
$ws.Range("F48").Value = -7
$ws.Range("F49").Value = -1
$ws.Range("F51").Value = -2
$ws.Range("F52").Value = -6
$ws.Range("F53").Value = 5
$ws.Range("F54").Value = 14
